$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking price cells so Excel keeps them as literal text
# (matches the source data which stores prices/volumes as text strings, not numbers)

$ws.Range('D2').Value = '27.973.69'
$ws.Range('E2').Value = '  +0.08%  '
$ws.Range('D3').Value = '1.855.74'
$ws.Range('E3').Value = '  -0.74%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  +0.35%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.76'
$ws.Range('E5').Value = '  -0.08%  '
$ws.Range('E6').Value = '  +0.30%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5073'
$ws.Range('E7').Value = '  +1.71%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3809'
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08237'
$ws.Range('E9').Value = '  -7.79%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.107'
$ws.Range('E10').Value = '  -1.11%  '
$ws.Range('E11').Value = '  +0.18%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.180'
$ws.Range('E12').Value = '  -2.06%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.863.19'
$ws.Range('E13').Value = '  +0.11%  '
$ws.Range('B14').Value = 'Solana'
$ws.Range('C14').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.44'
$ws.Range('E14').Value = '  -1.13%  '
$ws.Range('E15').Value = '  -0.75%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.003'
$ws.Range('E16').Value = '  +0.28%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001094'
$ws.Range('E17').Value = '  -0.56%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '90.32'
$ws.Range('E18').Value = '  -0.55%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06603'
$ws.Range('E19').Value = '  -0.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.67'
$ws.Range('E20').Value = '  -1.32%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.003'
$ws.Range('E21').Value = '  +0.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.999'
$ws.Range('E22').Value = '  -1.75%  '
$ws.Range('D23').Value = '27.987.20'
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.03'
$ws.Range('E24').Value = '  -4.24%  '
$ws.Range('E25').Value = '  -1.89%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.542'
$ws.Range('E26').Value = '  +1.02%  '
$ws.Range('D27').Value = '2.073.26'
$ws.Range('E27').Value = '  -0.15%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '157.99'
$ws.Range('E28').Value = '  +0.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.35'
$ws.Range('E29').Value = '  -1.68%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '124.16'
$ws.Range('E30').Value = '  -1.62%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1055'
$ws.Range('E31').Value = '  -0.29%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.032'
$ws.Range('E32').Value = '  -2.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.585'
$ws.Range('E33').Value = '  +0.12%  '
$ws.Range('E34').Value = '  +0.32%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.571'
$ws.Range('E35').Value = '  +2.55%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06513'
$ws.Range('E36').Value = '  -0.33%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02402'
$ws.Range('E37').Value = '  +0.26%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2165'
$ws.Range('E38').Value = '  -1.36%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.200'
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.238'
$ws.Range('E40').Value = '  -4.38%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6368'
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.17'
$ws.Range('E42').Value = '  -4.30%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.853'
$ws.Range('E43').Value = '  -0.91%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6040'
$ws.Range('E44').Value = '  +0.62%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.13'
$ws.Range('E45').Value = '  -0.73%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.281'
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.657'
$ws.Range('E47').Value = '  -0.38%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.972'
$ws.Range('E48').Value = '  -0.25%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.203'
$ws.Range('E49').Value = '  -1.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '120.53'
$ws.Range('E50').Value = '  -0.83%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '78.82'
$ws.Range('E51').Value = '  +0.95%  '
